$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Tamamlandı"
$ws.Cells.Item(1, 2).Value = "İL"
$ws.Cells.Item(1, 3).Value = "İLÇE"
$ws.Cells.Item(1, 4).Value = "BİRİM"
$ws.Cells.Item(1, 5).Value = "UYGULAMA"
$ws.Cells.Item(1, 6).Value = "GÖREVLİ PERSONELLER"
$ws.Cells.Item(1, 7).Value = "PARSEL SAYISI"
$ws.Cells.Item(1, 8).Value = "ALAN(Ha)"
$ws.Cells.Item(1, 9).Value = "İHALELİ/MÜDÜRLÜK"
$ws.Cells.Item(1, 10).Value = "İŞE BAŞLAMA/YER TESLİMİ"
$ws.Cells.Item(1, 11).Value = "İHALE BİTİŞ TARİHİ"
$ws.Cells.Item(1, 12).Value = "DURUMU"

# Row 2
$ws.Cells.Item(2, 1).Value = "HAYIR"
$ws.Cells.Item(2, 2).Value = "Adana"
$ws.Cells.Item(2, 3).Value = "Akdeniz"
$ws.Cells.Item(2, 4).Value = "ghghg"
$ws.Cells.Item(2, 5).Value = "GÜNCELLEME"
$ws.Cells.Item(2, 6).Value = "ESMEN TOKALI (K.Mühendisi), TEVFİK YILDIZ (Mühendis)"
$ws.Cells.Item(2, 7).Value = "'"
$ws.Cells.Item(2, 8).Value = "'4"
$ws.Cells.Item(2, 9).Value = "'7"
$ws.Cells.Item(2, 10).Value = "İhaleli"
$ws.Cells.Item(2, 11).Value = "'2025-02-01"
$ws.Cells.Item(2, 12).Value = "'2025-07-24"

# Row 3
$ws.Cells.Item(3, 1).Value = "HAYIR"
$ws.Cells.Item(3, 2).Value = "Adana"
$ws.Cells.Item(3, 3).Value = "Akdeniz"
$ws.Cells.Item(3, 4).Value = "hghmh"
$ws.Cells.Item(3, 5).Value = "2/B"
$ws.Cells.Item(3, 6).Value = "MEHMET AKGÜN KOLUKIRIK (Mühendis), AYHAN KARADAYI (K.Teknisyeni)"
$ws.Cells.Item(3, 7).Value = "'"
$ws.Cells.Item(3, 8).Value = "'4"
$ws.Cells.Item(3, 9).Value = "'12"
$ws.Cells.Item(3, 10).Value = "İhaleli"
$ws.Cells.Item(3, 11).Value = "'2025-01-01"
$ws.Cells.Item(3, 12).Value = "'2025-07-24"
